# Add rows 5-13 to sheet1: alternating literal dates / "+1" formulas in
# column A, and the "Shri Janardhana Swamy Temple" place name in column B,
# mirroring the existing row 3/4 pattern (date number format + wrapped,
# 60pt-tall place-name cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$placeName = "Shri Janardhana Swamy Temple"

# Column A: literal date in odd rows, "previous cell + 1" formula in even rows.
$ws.Range("A5").Value = 45736
$ws.Range("A6").Formula = "=A5+1"
$ws.Range("A7").Value = 45737
$ws.Range("A8").Formula = "=A7+1"
$ws.Range("A9").Value = 45738
$ws.Range("A10").Formula = "=A9+1"
$ws.Range("A11").Value = 45739
$ws.Range("A12").Formula = "=A11+1"
$ws.Range("A13").Value = 45740

# Column B: place name, same text for every new row.
for ($r = 5; $r -le 13; $r++) {
    $ws.Range("B$r").Value = $placeName
}

# Formatting: column A keeps the short-date number format used by A3/A4;
# column B keeps the wrap-text style used by B3, and every new row gets the
# same 60pt row height as row 3.
$ws.Range("A5:A13").NumberFormat = "m/d/yyyy"
$ws.Range("B5:B13").WrapText = $true

for ($r = 5; $r -le 13; $r++) {
    $ws.Rows.Item($r).RowHeight = 60
}

# Match the recorded post-edit selection.
$ws.Range("B5:B13").Select()
